$wb = $excel.ActiveWorkbook

# Sheet references
$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update status text "Ready for handoff" -> "In Translation".
# Find/Replace updates the shared string content itself (rather than
# creating a brand-new shared-string entry per cell), which is how a
# single status value used on all three sheets changes everywhere at once.
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("Ready for handoff", "In Translation", 1, 1, $false, $false, $false, $false) | Out-Null
}

# Shrink the now-narrower status columns to fit the shorter text:
#   Overview!E:F ("zh-cn" / "de-de" status columns)
#   zh-cn!C and de-de!C ("Status" column)
# Target stored width ~= 13.41 characters.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
